$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.433.16"
$ws.Range("E2").Value = "  +0.02%  "
$ws.Range("D3").Value = "1.941.42"
$ws.Range("E3").Value = "  -1.91%  "
$ws.Range("E4").Value = "  -0.25%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "242.28"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.27%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.603"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.30%  "
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "57.27"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.52%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.359"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.03%  "
$ws.Range("E10").Value = "  -2.86%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.103"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.47%  "
$ws.Range("D12").Value = "2.225.64"
$ws.Range("E12").Value = "  -2.03%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.31"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.61%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.810"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -5.16%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "13.46"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.78%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.14"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -5.66%  "
$ws.Range("D17").Value = "1.940.51"
$ws.Range("E17").Value = "  -1.21%  "
$ws.Range("D18").Value = "36.403.88"
$ws.Range("E18").Value = "  +0.27%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "69.14"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.76%  "
$ws.Range("D20").Value = "0.0₃0864"
$ws.Range("E20").Value = "  -4.75%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "227.77"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.95%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.99"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -5.10%  "
$ws.Range("E23").Value = "  -0.07%  "
$ws.Range("E24").Value = "  -6.50%  "
$ws.Range("E25").Value = "  -0.86%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.21"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -5.76%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "161.37"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.16%  "
$ws.Range("E28").Value = "  -0.04%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.21"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.63%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.117"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.11%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.10"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -7.32%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.56"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -6.43%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0618"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.43%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.17"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.68%  "
$ws.Range("E35").Value = "  -0.19%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.09"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.18%  "
$ws.Range("E37").Value = "  -1.02%  "
$ws.Range("E38").Value = "  -1.60%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.13"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +7.11%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0991"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.77%  "
$ws.Range("E41").Value = "  +0.15%  "
$ws.Range("E42").Value = "  -1.98%  "
$ws.Range("E43").Value = "  -4.83%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "15.66"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.15%  "
$ws.Range("D45").Value = "1.340.85"
$ws.Range("E45").Value = "  -1.68%  "
$ws.Range("E46").Value = "  -5.62%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "86.41"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.08%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.13"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.86%  "
$ws.Range("E49").Value = "  -0.26%  "
$ws.Range("D50").Value = "2.116.18"
$ws.Range("E50").Value = "  -2.02%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "43.07"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -4.88%  "
